# Updated symbol list on Mon Dec 12 20:30:48 UTC 2022 with GitHub Actions
#
# Applies the per-row "Price" (column D) and, where applicable, the
# "Volume(1h)" (column E) refreshes captured in the scraped diff.
#
# The Price column holds numeric-looking values that are stored as TEXT
# (t="inlineStr" in the original workbook), so a plain
# `$ws.Range(...).Value = "123"` would be auto-coerced to a Number by
# Excel. To preserve the original Text cell type we flip the cell to the
# Text number format before writing, then restore the cell's style so no
# residual formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Column D ("Price") updates
Set-TextValue "D2"  "275.05"
Set-TextValue "D3"  "21.17"
Set-TextValue "D4"  "6.256"
Set-TextValue "D5"  "0.06183"
Set-TextValue "D6"  "3.578"
Set-TextValue "D7"  "1.524"
Set-TextValue "D8"  "6.531"
Set-TextValue "D9"  "0.8228"
Set-TextValue "D11" "0.08263"
Set-TextValue "D13" "0.03163"
Set-TextValue "D14" "0.09143"
Set-TextValue "D16" "0.001613"
Set-TextValue "D17" "0.04681"
Set-TextValue "D18" "0.006260"
Set-TextValue "D20" "0.001069"
Set-TextValue "D21" "0.0001501"
Set-TextValue "D22" "3.725"
Set-TextValue "D23" "2.310"
Set-TextValue "D24" "0.01390"
Set-TextValue "D26" "0.1231"
Set-TextValue "D28" "0.0002738"
Set-TextValue "D40" "0.04738"
Set-TextValue "D41" "0.007041"
Set-TextValue "D42" "0.004454"
Set-TextValue "D44" "0.01149"
Set-TextValue "D45" "0.00006060"
Set-TextValue "D47" "0.7233"
Set-TextValue "D48" "0.001387"
Set-TextValue "D49" "0.00001901"
Set-TextValue "D50" "0.01241"

# Column E ("Volume(1h)") updates - plain text, no coercion risk
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
